# Apply MISSOURI_2022 sheet cleanup: rename headers, title-case connector words,
# correct tiny float roundoffs, and drop trailing footer/source rows 627-632.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header labels and title-case municipality/state connector words ---
$ws.Cells.Item(1, 1).Value2 = 'mx_state'
$ws.Cells.Item(1, 2).Value2 = 'mx_municipality'
$ws.Cells.Item(1, 3).Value2 = 'n_matriculas'
$ws.Cells.Item(1, 4).Value2 = 'pct_matriculas'
$ws.Cells.Item(4, 2).Value2 = 'Pabellón De Arteaga'
$ws.Cells.Item(5, 2).Value2 = 'Rincón De Romos'
$ws.Cells.Item(9, 2).Value2 = 'Playas De Rosarito'
$ws.Cells.Item(29, 2).Value2 = 'Marqués De Comillas'
$ws.Cells.Item(51, 2).Value2 = 'Hidalgo Del Parral'
$ws.Cells.Item(61, 2).Value2 = 'San Francisco Del Oro'
$ws.Cells.Item(71, 2).Value2 = 'San Juan De Sabinas'
$ws.Cells.Item(83, 1).Value2 = 'Ciudad De México'
$ws.Cells.Item(118, 1).Value2 = 'Estado De México'
$ws.Cells.Item(118, 2).Value2 = 'Acambay De Ruíz Castañeda'
$ws.Cells.Item(121, 2).Value2 = 'Almoloya De Juárez'
$ws.Cells.Item(123, 2).Value2 = 'Atizapán De Zaragoza'
$ws.Cells.Item(129, 2).Value2 = 'Ecatepec De Morelos'
$ws.Cells.Item(133, 2).Value2 = 'Ixtapan De La Sal'
$ws.Cells.Item(138, 2).Value2 = 'Naucalpan De Juárez'
$ws.Cells.Item(142, 2).Value2 = 'San Felipe Del Progreso'
$ws.Cells.Item(151, 2).Value2 = 'Tlalnepantla De Baz'
$ws.Cells.Item(155, 2).Value2 = 'Villa Del Carbón'
$ws.Cells.Item(162, 2).Value2 = 'Apaseo El Alto'
$ws.Cells.Item(165, 2).Value2 = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Cells.Item(168, 2).Value2 = 'Jaral Del Progreso'
$ws.Cells.Item(177, 2).Value2 = 'San Francisco Del Rincón'
$ws.Cells.Item(179, 2).Value2 = 'Santa Cruz De Juventino Rosas'
$ws.Cells.Item(180, 2).Value2 = 'Silao De La Victoria'
$ws.Cells.Item(187, 2).Value2 = 'Acapulco De Juárez'
$ws.Cells.Item(189, 2).Value2 = 'Ajuchitlán Del Progreso'
$ws.Cells.Item(190, 2).Value2 = 'Alcozauca De Guerrero'
$ws.Cells.Item(193, 2).Value2 = 'Atenango Del Río'
$ws.Cells.Item(194, 2).Value2 = 'Atoyac De Álvarez'
$ws.Cells.Item(195, 2).Value2 = 'Ayutla De Los Libres'
$ws.Cells.Item(197, 2).Value2 = 'Buenavista De Cuéllar'
$ws.Cells.Item(198, 2).Value2 = 'Chilapa De Álvarez'
$ws.Cells.Item(199, 2).Value2 = 'Chilpancingo De Los Bravo'
$ws.Cells.Item(202, 2).Value2 = 'Coyuca De Benítez'
$ws.Cells.Item(203, 2).Value2 = 'Coyuca De Catalán'
$ws.Cells.Item(205, 2).Value2 = 'Cutzamala De Pinzón'
$ws.Cells.Item(209, 2).Value2 = 'Huitzuco De Los Figueroa'
$ws.Cells.Item(217, 2).Value2 = 'Técpan De Galeana'
$ws.Cells.Item(218, 2).Value2 = 'Tlapa De Comonfort'
$ws.Cells.Item(228, 2).Value2 = 'Cuautepec De Hinojosa'
$ws.Cells.Item(233, 2).Value2 = 'Mixquiahuala De Juárez'
$ws.Cells.Item(234, 2).Value2 = 'Molango De Escamilla'
$ws.Cells.Item(235, 2).Value2 = 'Omitlán De Juárez'
$ws.Cells.Item(236, 2).Value2 = 'Pachuca De Soto'
$ws.Cells.Item(237, 2).Value2 = 'Progreso De Obregón'
$ws.Cells.Item(240, 2).Value2 = 'Tepehuacán De Guerrero'
$ws.Cells.Item(241, 2).Value2 = 'Tepeji Del Río De Ocampo'
$ws.Cells.Item(243, 2).Value2 = 'Tulancingo De Bravo'
$ws.Cells.Item(244, 2).Value2 = 'Zacualtipán De Ángeles'
$ws.Cells.Item(247, 2).Value2 = 'Ahualulco De Mercado'
$ws.Cells.Item(249, 2).Value2 = 'Atotonilco El Alto'
$ws.Cells.Item(250, 2).Value2 = 'Autlán De Navarro'
$ws.Cells.Item(252, 2).Value2 = 'Cuautitlán De García Barragán'
$ws.Cells.Item(256, 2).Value2 = 'Jilotlán De Los Dolores'
$ws.Cells.Item(259, 2).Value2 = 'La Manzanilla De La Paz'
$ws.Cells.Item(260, 2).Value2 = 'Lagos De Moreno'
$ws.Cells.Item(264, 2).Value2 = 'Ojuelos De Jalisco'
$ws.Cells.Item(266, 2).Value2 = 'San Juan De Los Lagos'
$ws.Cells.Item(267, 2).Value2 = 'San Miguel El Alto'
$ws.Cells.Item(268, 2).Value2 = 'San Sebastián Del Oeste'
$ws.Cells.Item(270, 2).Value2 = 'Tamazula De Gordiano'
$ws.Cells.Item(272, 2).Value2 = 'Tepatitlán De Morelos'
$ws.Cells.Item(273, 2).Value2 = 'Tlajomulco De Zúñiga'
$ws.Cells.Item(276, 2).Value2 = 'Unión De San Antonio'
$ws.Cells.Item(280, 2).Value2 = 'Zapotlán El Grande'
$ws.Cells.Item(289, 2).Value2 = 'Coalcomán De Vázquez Pallares'
$ws.Cells.Item(318, 2).Value2 = 'Tiquicheo De Nicolás Romero'
$ws.Cells.Item(340, 2).Value2 = 'Amatlán De Cañas'
$ws.Cells.Item(341, 2).Value2 = 'Santa María Del Oro'
$ws.Cells.Item(350, 2).Value2 = 'Montemorelos'
$ws.Cells.Item(353, 2).Value2 = 'Acatlán De Pérez Figueroa'
$ws.Cells.Item(354, 2).Value2 = 'Ayoquezco De Aldama'
$ws.Cells.Item(356, 2).Value2 = 'Chalcatongo De Hidalgo'
$ws.Cells.Item(357, 2).Value2 = 'Coicoyán De Las Flores'
$ws.Cells.Item(358, 2).Value2 = 'Cuilápam De Guerrero'
$ws.Cells.Item(359, 2).Value2 = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Cells.Item(364, 2).Value2 = 'Miahuatlán De Porfirio Díaz'
$ws.Cells.Item(365, 2).Value2 = 'Oaxaca De Juárez'
$ws.Cells.Item(401, 2).Value2 = 'Santo Domingo De Morelos'
$ws.Cells.Item(403, 2).Value2 = 'Tataltepec De Valdés'
$ws.Cells.Item(404, 2).Value2 = 'Tezoatlán De Segura Y Luna'
$ws.Cells.Item(405, 2).Value2 = 'Tlacolula De Matamoros'
$ws.Cells.Item(406, 2).Value2 = 'Villa De Tututepec De Melchor Ocampo'
$ws.Cells.Item(407, 2).Value2 = 'Villa Sola De Vega'
$ws.Cells.Item(408, 2).Value2 = 'Zimatlán De Álvarez'
$ws.Cells.Item(415, 2).Value2 = 'Chalchicomula De Sesma'
$ws.Cells.Item(419, 2).Value2 = 'Cuayuca De Andrade'
$ws.Cells.Item(422, 2).Value2 = 'Izúcar De Matamoros'
$ws.Cells.Item(431, 2).Value2 = 'San Salvador El Verde'
$ws.Cells.Item(437, 2).Value2 = 'Tepatlaxco De Hidalgo'
$ws.Cells.Item(439, 2).Value2 = 'Tlacotepec De Benito Juárez'
$ws.Cells.Item(451, 2).Value2 = 'Jalpan De Serra'
$ws.Cells.Item(454, 2).Value2 = 'San Juan Del Río'
$ws.Cells.Item(472, 2).Value2 = 'Villa De Ramos'
$ws.Cells.Item(495, 2).Value2 = 'Jalpa De Méndez'
$ws.Cells.Item(513, 2).Value2 = 'Soto La Marina'
$ws.Cells.Item(520, 2).Value2 = 'Acuamanala De Miguel Hidalgo'
$ws.Cells.Item(524, 2).Value2 = 'Contla De Juan Cuamatzi'
$ws.Cells.Item(526, 2).Value2 = 'Ixtacuixtla De Mariano Matamoros'
$ws.Cells.Item(528, 2).Value2 = 'Nanacamilpa De Mariano Arista'
$ws.Cells.Item(530, 2).Value2 = 'San Pablo Del Monte'
$ws.Cells.Item(533, 2).Value2 = 'Tetla De La Solidaridad'
$ws.Cells.Item(543, 2).Value2 = 'Alto Lucero De Gutiérrez Barrios'
$ws.Cells.Item(553, 2).Value2 = 'Cosamaloapan De Carpio'
$ws.Cells.Item(559, 2).Value2 = 'Hueyapan De Ocampo'
$ws.Cells.Item(560, 2).Value2 = 'Ignacio De La Llave'
$ws.Cells.Item(567, 2).Value2 = 'Lerdo De Tejada'
$ws.Cells.Item(569, 2).Value2 = 'Martínez De La Torre'
$ws.Cells.Item(572, 2).Value2 = 'Nanchital De Lázaro Cárdenas Del Río'
$ws.Cells.Item(577, 2).Value2 = 'Paso De Ovejas'
$ws.Cells.Item(579, 2).Value2 = 'Poza Rica De Hidalgo'
$ws.Cells.Item(584, 2).Value2 = 'Soledad De Doblado'
$ws.Cells.Item(606, 2).Value2 = 'Cañitas De Felipe Pescador'
$ws.Cells.Item(622, 2).Value2 = 'Villa De Cos'

# --- Correct minor floating point roundoffs in percentage column ---
$ws.Cells.Item(51, 4).Value2 = 0.009519038076152304
$ws.Cells.Item(158, 4).Value2 = 0.009018036072144287
$ws.Cells.Item(190, 4).Value2 = 0.009018036072144287
$ws.Cells.Item(210, 4).Value2 = 0.009018036072144287

# --- Remove trailing footnote/source rows (now rows 627-632 are gone; data ends at 626) ---
$ws.Rows("627:632").Delete()

# --- Ensure used range / dimension reflects the trimmed data (A1:D626) ---
Write-Host "Done"
